$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset any existing AutoFilter so the new one can size to the full data range
$ws.AutoFilterMode = $false

# Apply AutoFilter on column B ("Country") for "Bhutan" over the full data range
$ws.Range("A1:F193").AutoFilter(2, @("Bhutan"), 7)

# Keep the hidden _xlnm._FilterDatabase defined name in sync with the new range
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$F`$193"

# Move the active selection, matching the author's last cursor position
$ws.Range("E52").Select()
